$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-26 Monday", 2) | Out-Null

# Update each math-expression cell in the table, in row-major order,
# using positional cell addressing (some old values repeat, so a global
# Find/Replace would be ambiguous).
$newValues = @(
    "99-29=70",
    "47-11=36",
    "74-17=57",
    "16-10=6",
    "36+25=61",
    "77-72=5",
    "77-63=14",
    "35-4=31",
    "84+2=86",
    "4+55=59",
    "41-39=2",
    "30-26=4",
    "11-0=11",
    "3+71=74",
    "29+21=50",
    "76-3=73",
    "77-15=62",
    "39+19=58",
    "2+56=58",
    "47+1=48",
    "17-10=7",
    "45+9=54",
    "39+18=57",
    "82-72=10",
    "60+39=99",
    "25+41=66",
    "54-36=18",
    "4+70=74",
    "36+34=70",
    "89-6=83",
    "15+45=60",
    "12+7=19",
    "86-74=12",
    "16-14=2",
    "93-0=93",
    "9+32=41",
    "61-58=3",
    "62-40=22",
    "81-52=29",
    "75+18=93",
    "92-73=19",
    "74+18=92",
    "59-58=1",
    "54-20=34",
    "90-19=71",
    "50-42=8",
    "34+44=78",
    "53+31=84",
    "75-71=4",
    "8+15=23",
    "71-35=36",
    "30+61=91",
    "33-7=26",
    "63-10=53",
    "41-0=41",
    "24+47=71",
    "7+60=67",
    "6+64=70",
    "93-27=66",
    "56+27=83",
    "82-35=47",
    "85-68=17",
    "45-0=45",
    "43+36=79",
    "54-33=21",
    "36-11=25",
    "13+80=93",
    "50-9=41",
    "19-16=3",
    "71-65=6",
    "72-32=40",
    "79-7=72",
    "25+4=29",
    "7+39=46",
    "61+20=81",
    "90+8=98",
    "1+63=64",
    "42-0=42",
    "94-12=82",
    "13+81=94",
    "60-4=56",
    "33-24=9",
    "32-13=19",
    "71-32=39",
    "5+23=28",
    "35+45=80",
    "54-13=41",
    "49+44=93",
    "47+1=48",
    "23-6=17",
    "37-4=33",
    "13+8=21",
    "33+47=80",
    "80-25=55",
    "75-46=29",
    "18+23=41",
    "86+5=91",
    "95-73=22",
    "80-3=77",
    "41+46=87"
)

$table = $d.Tables.Item(1)
$cols = 5
$i = 0
foreach ($row in 1..$table.Rows.Count) {
    foreach ($col in 1..$cols) {
        $cell = $table.Cell($row, $col)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}

Write-Output ("Updated " + $i + " cells")